# "the last version 12 dec"
#
# - add one more CV link as a new last row on the "black_list" sheet
# - drop the explicit (wrapped-text) row heights on the earlier rows so
#   they go back to auto/default height
# - leave the view scrolled down near the new row, with it selected

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("black_list")

# Rows 13:41 previously carried explicit ht="..." (60 / 33.75 / 56.25 / 67.5 / 45)
# from when the cells were manually sized; AutoFit removes the hard-coded
# height so the row goes back to the sheet's default/auto height.
$ws.Range("A13:A41").Rows.AutoFit()

# New last entry.
$ws.Range("A42").Value = "https://rabota.ua/cv/3787170"

# Scroll the window down so row 28 is at the top and select the next
# (still empty) row, matching where editing left off.
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A43").Select()
